$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row just above the current row 29. This pushes the
# existing rows 29-70 down to 30-71 (preserving all of their data/format),
# and leaves an empty row 29 ready to be populated with the new record.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly record.
$ws.Range("A29").Value = 11
$ws.Range("B29").Value = "Vega Monumental Concepción"
$ws.Range("C29").Value = "Bíobío"
$ws.Range("D29").Value = 44965
$ws.Range("E29").Value = 8
$ws.Range("F29").Value = 100112031
$ws.Range("G29").Value = "Poroto verde"
$ws.Range("H29").Value = "Magnum"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 100
$ws.Range("K29").Value = 22000
$ws.Range("L29").Value = 23000
$ws.Range("M29").Value = 22500
$ws.Range("N29").Value = "`$/saco 25 kilos"
$ws.Range("O29").Value = "Región de O'Higgins"
$ws.Range("P29").Value = 900
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"
